$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 54 (shifts existing rows 54-93 down to 55-94),
# matching the diff's net effect: a new Alcachofa "Argentina(o)" price
# record was added for Terminal Hortofrutícola Agro Chillán.
$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value = 7
$ws.Range("B54").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C54").Value = "Ñuble"
$ws.Range("D54").Value = 45096
$ws.Range("E54").Value = 16
$ws.Range("F54").Value = 100112013
$ws.Range("G54").Value = "Alcachofa"
$ws.Range("H54").Value = "Argentina(o)"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 60
$ws.Range("K54").Value = 17000
$ws.Range("L54").Value = 17000
$ws.Range("M54").Value = 17000
$ws.Range("N54").Value = "`$/caja 50 unidades"
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 340
$ws.Range("Q54").Value = 50
$ws.Range("R54").Value = "Hortaliza"
